$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 37, pushing the existing rows 37-46 down to 38-47
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new weekly price record
$ws.Cells.Item(37, 1).Value  = 1
$ws.Cells.Item(37, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(37, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(37, 4).Value  = 44855
$ws.Cells.Item(37, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(37, 5).Value  = 15
$ws.Cells.Item(37, 6).Value  = 100112052
$ws.Cells.Item(37, 7).Value  = "Albahaca"
$ws.Cells.Item(37, 8).Value  = "Sin especificar"
$ws.Cells.Item(37, 9).Value  = "Primera"
$ws.Cells.Item(37, 10).Value = 600
$ws.Cells.Item(37, 11).Value = 3000
$ws.Cells.Item(37, 12).Value = 3500
$ws.Cells.Item(37, 13).Value = 3167
$ws.Cells.Item(37, 14).Value = "$/atado"
$ws.Cells.Item(37, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(37, 16).Value = 3167
$ws.Cells.Item(37, 17).Value = 1
$ws.Cells.Item(37, 18).Value = "Hortaliza"
